# Update the "PPRiFUfIIaIoE" sheet's header cell (B1) to add the
# "(dimensionless)" unit suffix, and make it wrap so the taller header
# row matches the published workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PPRiFUfIIaIoE")

$ws.Range("B1").Value = "Pot Perc Red in Fuel Use (dimensionless)"
$ws.Range("B1").WrapText = $true
$ws.Rows.Item(1).RowHeight = 28.5

# Restore the original author's view state: cell B1 was left selected on
# this sheet (even though the "About" sheet is the one active on save).
[void]$ws.Range("B1").Select()
[void]$wb.Worksheets.Item("About").Activate()
